$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header/index-column style (bold, centered, thin border) used by
# the existing header row (e.g. H1) on the two new header cells.
$hdr = $ws.Range("I1:J1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

# Data for columns I (I0) and J (IF) for rows 2..68
$iVals = @(8,4,9,9,9,8,8,9,9,8,9,9,9,9,9,8,9,9,10,9,9,8,9,10,9,9,9,8,8,9,10,9,10,9,8,9,9,9,8,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,10,10,9,9,9,9,9,9,9,9,9,4,4)
$jVals = @(9,4,9,9,9,9,9,9,9,8,9,9,9,9,9,9,9,10,10,9,9,9,10,10,9,9,9,9,9,9,10,9,10,9,8,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,10,10,9,9,9,9,9,9,9,9,9,4,4)

for ($i = 0; $i -lt $iVals.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 9).Value = $iVals[$i]
    $ws.Cells.Item($r, 10).Value = $jVals[$i]
}
